$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Flip F-column "Status" for several rotations.py comparison rows
#    from "Neutral" (orange fill) to "Win" (green fill), matching the
#    fill/style already used by existing "Win" rows (e.g. F2).
$winRows = @(67,69,71,73,75,77,83,91,97,99)
foreach ($r in $winRows) {
    $cell = $ws.Range("F$r")
    $cell.Interior.Color = $ws.Range("F2").Interior.Color
    $cell.Value = "Win"
}

# 2. Update the "Scalpel Wins:" total on the summary row (14 -> 24)
$ws.Range("F125").Value = 24

# 3. Add a new trailing summary row (127) that carries the
#    "Accuracy over PyType" figure down one row, and repurpose row 126
#    to hold the new "Scalpel Accuracy:" metric.
$ws.Range("A127:D127").Value = ""
$ws.Range("A127:F127").Interior.Color = $ws.Range("A126").Interior.Color
$ws.Range("E127").Value = "Accuracy over PyType"
$ws.Range("F127").Value = 218.18

$ws.Range("C126").Value = "Scalpel Accuracy:"
$ws.Range("D126").Value = 1018.18
$ws.Range("E126").Value = ""
$ws.Range("F126").Value = ""
